$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("H3").Value = 0.8901191163352049
$ws.Range("I3").Value = 0.008220923490096744
$ws.Range("K3").Value = 344.3225806451613

$ws.Range("Q3").Value = 8
$ws.Range("R3").Value = 31
$ws.Range("S3").Value = 93
$ws.Range("T3").Value = 223
$ws.Range("U3").Value = 597

$ws.Range("V3").Value = 38044
$ws.Range("W3").Value = 38021
$ws.Range("X3").Value = 37959
$ws.Range("Y3").Value = 37829
$ws.Range("Z3").Value = 37455

$ws.Range("AF3").Value = 0.99979
$ws.Range("AG3").Value = 0.999185
$ws.Range("AH3").Value = 0.997556
$ws.Range("AI3").Value = 0.99414
$ws.Range("AJ3").Value = 0.984311
